# Apply "Trade #11" update across the workbook.

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.81
$wsSummary.Range("B4").Value = -0.19
$wsSummary.Range("B5").Value = -0.35
$wsSummary.Range("B6").Value = 11
$wsSummary.Range("B8").Value = 8
$wsSummary.Range("B9").Value = 27.27

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.81
$wsStatus.Range("D4").Value = 11
$wsStatus.Range("E4").Value = -0.19
$wsStatus.Range("F4").Value = -0.19
$wsStatus.Range("G4").Value = 27.27

# --- Append new trade row (#11) to "All Trades" and "MarketMaking" sheets ---
$newRow = @(
    11,
    "2026-02-17",
    "13:34:37",
    "MarketMaking",
    "DOWN",
    0.97,
    0.9,
    "CLOSED",
    -7.2165,
    -0.07000000000000001,
    99.81,
    0,
    0,
    0.6,
    "Normal spread capture: 19600 bps",
    "early_exit",
    0.1
)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Ensure the date-looking text value (column B) is kept as plain text,
    # matching the rest of the column, rather than being auto-converted to
    # an Excel date serial number.
    $ws.Range("B12").NumberFormat = "@"

    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $ws.Cells.Item(12, $i + 1).Value = $newRow[$i]
    }
}
